$wb = $excel.ActiveWorkbook

# Remember the currently active sheet so we can restore it at the end
# (selecting a range on another sheet implicitly activates that sheet).
$activeSheet = $wb.ActiveSheet

$ws = $wb.Worksheets.Item("BCS-BCS")

# Remove the IRA/EPA-regs-driven subsidy formulas in D2:M2 and D3:M3,
# replacing them with literal 0 values (same number format as the
# neighboring B:C columns).
$ws.Range("D2:M3").Value = 0
$ws.Range("D2:M3").NumberFormat = "0"

# Update the selection on this sheet to match the edited range.
$null = $ws.Range("C2:M3").Select()

# Restore the original active sheet/tab.
$null = $activeSheet.Activate()
